$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original table (Table1, A1:B7):
#   Row1: id_sector | name
#   Row2: 1 | Agriculture
#   Row3: 2 | Industry
#   Row4: 3 | Tertiary
#   Row5: 4 | Construction
#   Row6: 5 | Energy
#   Row7: 6 | Residential
#
# Target table (Table1, A1:B3) keeps only the Tertiary and Residential
# sectors used for the cooling/ventilation capex & opex calculation.

$table = $ws.ListObjects.Item("Table1")

# Replace data rows 2 and 3 with the sectors we keep (Tertiary, Residential)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Tertiary"
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "Residential"

# Remove the now unused rows (old rows 4-7: Construction, Energy, Residential dup, blank)
$ws.Range("A4:B7").EntireRow.Delete()

# Shrink the table to the new extent
$table.Resize($ws.Range("A1:B3"))

# Leave the selection where the author left it
$ws.Range("D11").Select()
